# Daily attendance processing - reorders the "Recorded By" (column G)
# values so that entries recorded alongside "dnasr281@gmail.com" list
# dnasr281's co-recorder first and dnasr281@gmail.com second, e.g.
#   "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#   "dnasr281@gmail.com, admin@admin.com" -> "admin@admin.com, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value()

    if ($val -ne $null -and $val.StartsWith("dnasr281@gmail.com, ")) {
        $parts = $val.Split(",")
        if ($parts.Count -eq 2) {
            $first = $parts[0].Trim()
            $second = $parts[1].Trim()
            $cell.Value = "$second, $first"
        }
    }
}
